$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows at position 166-167 (everything from old row 166
# downward shifts down by 2, ending at row 244 instead of row 242).
$ws.Rows("166:167").Insert()

# New row 166: Tomate, Larga vida, Primera, Vega Monumental Concepción
$ws.Range("A166").Value = 11
$ws.Range("B166").Value = "Vega Monumental Concepción"
$ws.Range("C166").Value = "Bíobío"
$ws.Range("D166").Value = 44460
$ws.Range("E166").Value = 8
$ws.Range("F166").Value = 100112020
$ws.Range("G166").Value = "Tomate"
$ws.Range("H166").Value = "Larga vida"
$ws.Range("I166").Value = "Primera"
$ws.Range("J166").Value = 600
$ws.Range("K166").Value = 15000
$ws.Range("L166").Value = 16000
$ws.Range("M166").Value = 15500
$ws.Range("N166").Value = "`$/bandeja 18 kilos"
$ws.Range("O166").Value = "Región de Arica y Parinacota"
$ws.Range("P166").Value = 861
$ws.Range("Q166").Value = 18
$ws.Range("R166").Value = "Hortaliza"

# New row 167: Tomate, Larga vida, Segunda, Vega Monumental Concepción
$ws.Range("A167").Value = 11
$ws.Range("B167").Value = "Vega Monumental Concepción"
$ws.Range("C167").Value = "Bíobío"
$ws.Range("D167").Value = 44460
$ws.Range("E167").Value = 8
$ws.Range("F167").Value = 100112020
$ws.Range("G167").Value = "Tomate"
$ws.Range("H167").Value = "Larga vida"
$ws.Range("I167").Value = "Segunda"
$ws.Range("J167").Value = 300
$ws.Range("K167").Value = 14000
$ws.Range("L167").Value = 14000
$ws.Range("M167").Value = 14000
$ws.Range("N167").Value = "`$/bandeja 18 kilos"
$ws.Range("O167").Value = "Región de Arica y Parinacota"
$ws.Range("P167").Value = 778
$ws.Range("Q167").Value = 18
$ws.Range("R167").Value = "Hortaliza"
